# Apply the cryptos list refresh (prices + 1h volume %) for Wed Mar  1 17:26:48 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $NewValue)
    $c = $ws.Range($CellRef)
    # Force the cell to be treated as Text so numeric-looking strings
    # (e.g. "303.11") are not silently converted to numbers by Excel,
    # matching the inline-string cells already used in this sheet.
    $c.NumberFormat = "@"
    $c.Value = $NewValue
    # Restore the default (unstyled) cell style so we do not leave
    # a stray number-format override behind on the cell itself.
    $c.Style = "Normal"
}

Set-TextValue 'D2' '23.728.65'
Set-TextValue 'D3' '1.658.93'
Set-TextValue 'E3' '  +1.08%  '
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  +0.20%  '
Set-TextValue 'D5' '1.000'
Set-TextValue 'E5' '  +0.08%  '
Set-TextValue 'D6' '303.11'
Set-TextValue 'E6' '  -0.31%  '
Set-TextValue 'D7' '0.3816'
Set-TextValue 'E7' '  +0.71%  '
Set-TextValue 'E8' '  -0.30%  '
Set-TextValue 'D9' '51.26'
Set-TextValue 'E9' '  -1.19%  '
Set-TextValue 'D10' '0.08194'
Set-TextValue 'E10' '  +0.07%  '
Set-TextValue 'D11' '1.229'
Set-TextValue 'E11' '  -0.54%  '
Set-TextValue 'D12' '1.001'
Set-TextValue 'E12' '  +0.05%  '
Set-TextValue 'D13' '22.54'
Set-TextValue 'E13' '  -0.25%  '
Set-TextValue 'D14' '6.471'
Set-TextValue 'E14' '  +0.08%  '
Set-TextValue 'D15' '7.405'
Set-TextValue 'E15' '  +0.32%  '
Set-TextValue 'D16' '0.00001228'
Set-TextValue 'E16' '  -1.00%  '
Set-TextValue 'D17' '1.653.79'
Set-TextValue 'E17' '  +1.49%  '
Set-TextValue 'D18' '97.93'
Set-TextValue 'E18' '  +2.75%  '
Set-TextValue 'D19' '0.07017'
Set-TextValue 'E19' '  +1.03%  '
Set-TextValue 'D20' '6.821'
Set-TextValue 'E20' '  +3.56%  '
Set-TextValue 'D21' '17.63'
Set-TextValue 'E21' '  +0.43%  '
Set-TextValue 'D22' '1.001'
Set-TextValue 'E22' '  +0.06%  '
Set-TextValue 'D23' '12.82'
Set-TextValue 'E23' '  +2.36%  '
Set-TextValue 'D24' '23.748.17'
Set-TextValue 'E24' '  +0.98%  '
Set-TextValue 'D25' '2.509'
Set-TextValue 'E25' '  +0.07%  '
Set-TextValue 'D26' '2.991'
Set-TextValue 'E26' '  -2.48%  '
Set-TextValue 'D27' '21.23'
Set-TextValue 'E27' '  +0.07%  '
Set-TextValue 'D28' '153.79'
Set-TextValue 'E28' '  +1.00%  '
Set-TextValue 'D29' '5.226'
Set-TextValue 'E29' '  -0.61%  '
Set-TextValue 'D30' '134.30'
Set-TextValue 'E30' '  +0.62%  '
Set-TextValue 'D31' '1.839.20'
Set-TextValue 'E31' '  +1.75%  '
Set-TextValue 'D32' '7.072'
Set-TextValue 'E32' '  +6.81%  '
Set-TextValue 'D33' '2.224'
Set-TextValue 'E33' '  +3.47%  '
Set-TextValue 'D34' '12.02'
Set-TextValue 'E34' '  +4.63%  '
Set-TextValue 'D35' '1.057'
Set-TextValue 'E35' '  -4.04%  '
Set-TextValue 'D36' '0.02814'
Set-TextValue 'E36' '  +1.77%  '
Set-TextValue 'D37' '0.2516'
Set-TextValue 'E37' '  +0.51%  '
Set-TextValue 'D38' '0.08803'
Set-TextValue 'E38' '  +0.25%  '
Set-TextValue 'D39' '6.086'
Set-TextValue 'E39' '  +1.02%  '
Set-TextValue 'D40' '0.07011'
Set-TextValue 'E40' '  -1.17%  '
Set-TextValue 'D41' '12.98'
Set-TextValue 'E41' '  +5.62%  '
Set-TextValue 'D42' '0.7005'
Set-TextValue 'E42' '  -0.94%  '
Set-TextValue 'D43' '1.335'
Set-TextValue 'E43' '  -1.25%  '
Set-TextValue 'E44' '  +2.85%  '
Set-TextValue 'D45' '0.6516'
Set-TextValue 'E45' '  -0.68%  '
Set-TextValue 'B46' 'NEARProtocol'
Set-TextValue 'C46' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D46' '2.313'
Set-TextValue 'E46' '  +1.18%  '
Set-TextValue 'B47' 'Frax'
Set-TextValue 'C47' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D47' '0.9995'
Set-TextValue 'E47' '  -0.03%  '
Set-TextValue 'D48' '3.969'
Set-TextValue 'E48' '  +0.01%  '
Set-TextValue 'D49' '0.07925'
Set-TextValue 'E49' '  -0.78%  '
Set-TextValue 'D50' '128.36'
Set-TextValue 'E50' '  -0.46%  '
Set-TextValue 'D51' '1.185'
Set-TextValue 'E51' '  -0.93%  '
